# Update "paises" workbook: refresh COVID stats + swap a few country rows
# back into (now) correct order, and bump the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 7 de Julio de 2020 a las 05:57"

# --- Brasil (row 5): refresh recuperados / casos activos --------------
$ws.Range("D5").Value = 1072229
$ws.Range("E5").Value = 488286

# --- Guatemala / Honduras swap (rows 55-56) ----------------------------
$ws.Range("A55").Value = "Honduras"
$ws.Range("B55").Value = 24665
$ws.Range("C55").Value = 722
$ws.Range("D55").Value = 2585
$ws.Range("E55").Value = 21424
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 17
$ws.Range("H55").Value = 656

$ws.Range("A56").Value = "Guatemala"
$ws.Range("B56").Value = 23972
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 3429
$ws.Range("E56").Value = 19562
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 981

# --- Malasia / Australia swap (rows 74-75) -----------------------------
$ws.Range("A74").Value = "Australia"
$ws.Range("B74").Value = 8755
$ws.Range("C74").Value = 169
$ws.Range("D74").Value = 7455
$ws.Range("E74").Value = 1194
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 106

$ws.Range("A75").Value = "Malasia"
$ws.Range("B75").Value = 8668
$ws.Range("C75").Value = 0
$ws.Range("D75").Value = 8476
$ws.Range("E75").Value = 71
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 121

# --- Haiti (row 84): refresh stats --------------------------------------
$ws.Range("B84").Value = 6371
$ws.Range("C84").Value = 38
$ws.Range("D84").Value = 1824
$ws.Range("E84").Value = 4434

# --- Mongolia (row 169): refresh stats ----------------------------------
$ws.Range("B169").Value = 225
$ws.Range("C169").Value = 5
$ws.Range("D169").Value = 194
$ws.Range("E169").Value = 31

# --- Nueva Caledonia / Fiyi swap (rows 203-204) -------------------------
$ws.Range("A203").Value = "Fiyi"
$ws.Range("B203").Value = 21
$ws.Range("C203").Value = 2
$ws.Range("D203").Value = 18
$ws.Range("E203").Value = 3

$ws.Range("A204").Value = "Nueva Caledonia"
$ws.Range("B204").Value = 21
$ws.Range("C204").Value = 0
$ws.Range("D204").Value = 21
$ws.Range("E204").Value = 0

# --- Groenlandia / Islas Malvinas swap (rows 209-210) -------------------
# (numeric values are identical for these two rows, only the names swap)
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"
